# Update the "Solar" capacity figures for 2021 and 2024 with revised
# upstream data (commit: "Incorporate updated data from upstream
# processes through 2024").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 2021 Solar capacity: 93.59999999999999 -> 90.3
$ws.Range("E23").Value = 90.3

# 2024 Solar capacity: 136.95 -> 146.95
$ws.Range("E26").Value = 146.95
